$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# STEP 1: Split paragraph 2 ("DyWorld has a RPG system, access it with its
# default key of NUMPAD 5.") into five separate runs:
#   " - DyWorld has a RPG system" | ". D" | "efault key" | ":" | " NUMPAD 5."
#
# Technique: shrink the paragraph to the first run's text, then repeatedly
# insert a fresh paragraph right after it holding the next run's text, and
# merge it back by deleting the paragraph mark that separates them. Because
# the merge happens by deleting a paragraph boundary (not by appending text
# into an existing run), the two pieces of text stay as distinct <w:r> runs
# even though their formatting is identical.
# ---------------------------------------------------------------------------

$p = $d.Paragraphs.Item(2)
$p.Range.Text = " - DyWorld has a RPG system"

$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(3)
$np.Range.Text = ". D"
$p = $d.Paragraphs.Item(2)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(3)
$np.Range.Text = "efault key"
$p = $d.Paragraphs.Item(2)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(3)
$np.Range.Text = ":"
$p = $d.Paragraphs.Item(2)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(3)
$np.Range.Text = " NUMPAD 5."
$p = $d.Paragraphs.Item(2)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# ---------------------------------------------------------------------------
# STEP 2: Remove paragraph 3 ("Trees and Stones are the only thing the
# player can mine manually. ...") entirely -- that whole tip is gone.
# ---------------------------------------------------------------------------

$p = $d.Paragraphs.Item(3)
$p.Range.Delete()

# ---------------------------------------------------------------------------
# STEP 3: Remove the old "Use the Story menu ... NUMPAD 6." paragraph from
# its original spot (now paragraph 4, right after "player does NOT heal").
# It gets re-created right after the DyWorld tip below.
# ---------------------------------------------------------------------------

$p = $d.Paragraphs.Item(4)
$p.Range.Delete()

# ---------------------------------------------------------------------------
# STEP 4: Insert three new paragraphs right after paragraph 2 (DyWorld):
#   - "Use the Story menu to progress into the game. ... Default key: NUMPAD 6."
#   - "Want to reread the story? Use the log gui. Default key: NUMPAD 7"
#   - " - " + "Most Deposits can NOT be mined by the player! Use mining drills!" (2 runs)
# ---------------------------------------------------------------------------

$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(3)
$np.Range.Text = " - Use the Story menu to progress into the game. This unlocks recipes, technologies and a fancy story. Default key: NUMPAD 6."

$p = $d.Paragraphs.Item(3)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(4)
$np.Range.Text = " - Want to reread the story? Use the log gui. Default key: NUMPAD 7"

$p = $d.Paragraphs.Item(4)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(5)
$np.Range.Text = " - "

$p = $d.Paragraphs.Item(5)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(6)
$np.Range.Text = "Most Deposits can NOT be mined by the player! Use mining drills!"
$p = $d.Paragraphs.Item(5)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# ---------------------------------------------------------------------------
# STEP 5: Split the "Tier 0 (start of the game) is designed to be hard! ..."
# paragraph into three runs: " - " | "Act 1" | " (start of the game) is
# designed to be hard! You might die sometimes."
# ---------------------------------------------------------------------------

$p = $d.Paragraphs.Item(9)
$p.Range.Text = " - "

$p = $d.Paragraphs.Item(9)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(10)
$np.Range.Text = "Act 1"
$p = $d.Paragraphs.Item(9)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

$p = $d.Paragraphs.Item(9)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs.Item(10)
$np.Range.Text = " (start of the game) is designed to be hard! You might die sometimes."
$p = $d.Paragraphs.Item(9)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# ---------------------------------------------------------------------------
# STEP 6: Append a brand new tip paragraph at the very end of the document.
# ---------------------------------------------------------------------------

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Range.Text = " - Be careful with armor! They are not infinite, and especially gridded armor can be expensive with losing your added items"

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
